$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) is unchanged ("id", "Best Sports Book(HT)", "Favorite",
# "Favored Team Spread", "Underdog Team", "Underdog Team Spread",
# "Home Arena Belongs To").

# Clear out only the old data rows (2-7) so the new "Monday Slate" data can
# be written in without leaving now-unused shared strings (old team/book
# names) behind in the workbook.
$ws.Range("A2:G7").Clear()

# New data for rows 2-10 (id, favored sportsbook, favorite team,
# favored spread, underdog team, underdog spread,
# "home arena belongs to" team).
$colA = @(0, 1, 2, 3, 4, 5, 6, 7, 8)
$colB = @("FanDuel", "FanDuel", "FanDuel", "Unibet", "WynnBET", "William Hill (US)", "TwinSpires", "BetRivers", "LowVig.ag")
$colC = @("Charlotte Hornets", "Cleveland Cavaliers", "Detroit Pistons", "Miami Heat", "Brooklyn Nets", "Boston Celtics", "Chicago Bulls", "Houston Rockets", "Dallas Mavericks")
$colD = @(-6, -6, -7.5, -1, -1, -13.5, -4, 5, -3)
$colE = @("New Orleans Pelicans", "Los Angeles Lakers", "Portland Trail Blazers", "Philadelphia 76ers", "Utah Jazz", "Oklahoma City Thunder", "Toronto Raptors", "Washington Wizards", "Minnesota Timberwolves")
$colF = @(6, 6, 7.5, 1, 1, 13.5, 4, -5, 3)
$colG = @("Charlotte Hornets", "Cleveland Cavaliers", "Detroit Pistons", "Philadelphia 76ers", "Brooklyn Nets", "Oklahoma City Thunder", "Chicago Bulls", "Houston Rockets", "Dallas Mavericks")

# Write column-by-column (B, then C, then D, then E, then F, then G, then A)
# so new shared-string table entries are created in the same order the
# original export produced them.
for ($i = 0; $i -lt $colB.Count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $colB[$i]
}
for ($i = 0; $i -lt $colC.Count; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $colC[$i]
}
for ($i = 0; $i -lt $colD.Count; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $colD[$i]
}
for ($i = 0; $i -lt $colE.Count; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $colE[$i]
}
for ($i = 0; $i -lt $colF.Count; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $colF[$i]
}
for ($i = 0; $i -lt $colG.Count; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $colG[$i]
}
for ($i = 0; $i -lt $colA.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $colA[$i]
}

# Re-apply the bordered/centered "id column" formatting (originally style
# index 1, also used by the A1 header cell, which the row clear above did
# not touch) to the full new id range A2:A10, by copying the format from
# A1 - this reuses the existing style instead of creating new, unused
# style entries.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A2:A10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
